$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to Text
# storage (matching the source workbook, where the Price column is stored
# as inline strings, not numbers) by temporarily switching the cell to a
# text number format, then restoring the default ("Normal") style so no
# stray formatting is left behind.

$ws.Range("D2").Value = '87.397.29'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '3.157.56'
$ws.Range("E3").Value = '  -5.15%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '605.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.91%  '
$ws.Range("E7").Value = '  -5.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.666'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.88%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '3.152.51'
$ws.Range("E10").Value = '  -5.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.531'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -11.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.175'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -11.74%  '
$ws.Range("D14").Value = '3.740.06'
$ws.Range("E14").Value = '  -5.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '87.261.06'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.68%  '
$ws.Range("D18").Value = '3.142.76'
$ws.Range("E18").Value = '  -6.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '411.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.52%  '
$ws.Range("D26").Value = '3.325.44'
$ws.Range("E26").Value = '  -5.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '72.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000130'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.159'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '539.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("E33").Value = '  -9.77%  '
$ws.Range("E34").Value = '  -12.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.65'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.30%  '
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -7.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.366'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -12.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '171.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.92%  '
$ws.Range("E48").Value = '  +2.52%  '
$ws.Range("E49").Value = '  -12.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.61%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.80%  '
